$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 and J1, copying the format (style) used by the existing H1 header
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 values (column I), for rows 2..35
$i0 = @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,3,1,1,5)
# IF values (column J), for rows 2..35
$if = @(5,6,6,7,5,6,6,6,6,3,6,3,6,7,6,7,6,6,5,6,6,5,3,7,6,7,6,6,6,6,5,5,5,7)

for ($n = 0; $n -lt 34; $n++) {
    $row = 2 + $n
    $ws.Cells.Item($row, 9).Value = $i0[$n]
    $ws.Cells.Item($row, 10).Value = $if[$n]
}
